# NATMI ligand-receptor pair table update ("following Dr Hou advice"):
# the sending-cluster set is expanded from {FAPs, sCs} to {ECs, FAPs, sCs},
# which changes every edge's stats (ligand/receptor expression, specificity,
# weights) for the Ntn4 -> Unc5a pair against all four target clusters
# (ECs, FAPs, M2, sCs). Re-derive rows 2-13 (12 rows = 3 senders x 4 targets)
# from the new computation; row 1 (header) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..T per row:
# A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E=Ligand-expressing cells, F=Ligand detection rate,
# G=Ligand average expression value, H=Ligand total expression value,
# I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
# K=Receptor-expressing cells, L=Receptor detection rate,
# M=Receptor average expression value, N=Receptor total expression value,
# O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity
$rows = @(
    ,@("ECs", "Ntn4", "Unc5a", "ECs", 2, 0.6666666666666666, 1.664730333333334, 4.994191000000001, 0.05028984785586801, 0.05028984785586801, 3, 1, 2.322965, 6.968895, 0.2900587755589272, 0.2900587755589272, 3.867110298771667, 34.803992688945, 0.01458701169211782, 0.01458701169211782)
    ,@("ECs", "Ntn4", "Unc5a", "FAPs", 2, 0.6666666666666666, 1.664730333333334, 4.994191000000001, 0.05028984785586801, 0.05028984785586801, 3, 1, 1.626459, 4.879377, 0.2030890289077955, 0.2030890289077955, 2.707615633223, 24.368540699007, 0.01021331636496902, 0.01021331636496902)
    ,@("ECs", "Ntn4", "Unc5a", "M2", 2, 0.6666666666666666, 1.664730333333334, 4.994191000000001, 0.05028984785586801, 0.05028984785586801, 3, 1, 0.904191, 2.712573, 0.1129024907096758, 0.1129024907096758, 1.505234184827, 13.547107663443, 0.005677849080338146, 0.005677849080338145)
    ,@("ECs", "Ntn4", "Unc5a", "sCs", 2, 0.6666666666666666, 1.664730333333334, 4.994191000000001, 0.05028984785586801, 0.05028984785586801, 3, 1, 3.154986, 9.464958, 0.3939497048236015, 0.3939497048236015, 5.252200895442001, 47.269808058978, 0.01981167071844303, 0.01981167071844303)
    ,@("FAPs", "Ntn4", "Unc5a", "ECs", 3, 1, 12.206517, 36.619551, 0.3687467396301421, 0.3687467396301421, 3, 1, 2.322965, 6.968895, 0.2900587755589272, 0.2900587755589272, 28.355311762905, 255.197805866145, 0.1069582277884656, 0.1069582277884656)
    ,@("FAPs", "Ntn4", "Unc5a", "FAPs", 3, 1, 12.206517, 36.619551, 0.3687467396301421, 0.3687467396301421, 3, 1, 1.626459, 4.879377, 0.2030890289077955, 0.2030890289077955, 19.853399433303, 178.680594899727, 0.07488841726440128, 0.07488841726440128)
    ,@("FAPs", "Ntn4", "Unc5a", "M2", 3, 1, 12.206517, 36.619551, 0.3687467396301421, 0.3687467396301421, 3, 1, 0.904191, 2.712573, 0.1129024907096758, 0.1129024907096758, 11.037022812747, 99.333205314723, 0.04163242534531535, 0.04163242534531535)
    ,@("FAPs", "Ntn4", "Unc5a", "sCs", 3, 1, 12.206517, 36.619551, 0.3687467396301421, 0.3687467396301421, 3, 1, 3.154986, 9.464958, 0.3939497048236015, 0.3939497048236015, 38.51139024376199, 346.602512193858, 0.1452676692319599, 0.1452676692319599)
    ,@("sCs", "Ntn4", "Unc5a", "ECs", 3, 1, 19.23146433333334, 57.69439300000001, 0.58096341251399, 0.58096341251399, 3, 1, 2.322965, 6.968895, 0.2900587755589272, 0.2900587755589272, 44.67401854508167, 402.066166905735, 0.1685135360783439, 0.1685135360783439)
    ,@("sCs", "Ntn4", "Unc5a", "FAPs", 3, 1, 19.23146433333334, 57.69439300000001, 0.58096341251399, 0.58096341251399, 3, 1, 1.626459, 4.879377, 0.2030890289077955, 0.2030890289077955, 31.279188248129, 281.512694233161, 0.1179872952784252, 0.1179872952784252)
    ,@("sCs", "Ntn4", "Unc5a", "M2", 3, 1, 19.23146433333334, 57.69439300000001, 0.58096341251399, 0.58096341251399, 3, 1, 0.904191, 2.712573, 0.1129024907096758, 0.1129024907096758, 17.388916967021, 156.500252703189, 0.06559221628402229, 0.06559221628402227)
    ,@("sCs", "Ntn4", "Unc5a", "sCs", 3, 1, 19.23146433333334, 57.69439300000001, 0.58096341251399, 0.58096341251399, 3, 1, 3.154986, 9.464958, 0.3939497048236015, 0.3939497048236015, 60.675000731166, 546.075006580494, 0.2288703648731986, 0.2288703648731986)
)

# Write cell-by-cell (Range.Value2 array assignment isn't reliably
# marshalled in this host), starting at row 2 (row 1 is the header).
# This both overwrites the original 8 data rows (2-9) and appends the
# 4 new ones (10-13), extending the sheet dimension to A1:T13.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value2 = $row[$j]
    }
}

Write-Host "Wrote $($rows.Count) data rows (rows 2..$($rows.Count + 1))"
